$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.009.38'
$ws.Range('E2').Value = '  +3.24%  '
$ws.Range('D3').Value = '3.029.02'
$ws.Range('E3').Value = '  +1.79%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.51'
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('E6').Value = '  +7.69%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.026.22'
$ws.Range('E8').Value = '  +1.79%  '
$ws.Range('E9').Value = '  +0.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.94'
$ws.Range('E10').Value = '  +16.07%  '
$ws.Range('E11').Value = '  +2.85%  '
$ws.Range('E12').Value = '  +2.63%  '
$ws.Range('E13').Value = '  +3.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.52'
$ws.Range('E14').Value = '  +4.38%  '
$ws.Range('E15').Value = '  -0.40%  '
$ws.Range('D16').Value = '3.538.91'
$ws.Range('E16').Value = '  +2.15%  '
$ws.Range('E17').Value = '  +3.70%  '
$ws.Range('D18').Value = '62.946.96'
$ws.Range('E18').Value = '  +2.96%  '
$ws.Range('D19').Value = '3.031.16'
$ws.Range('E19').Value = '  +2.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '451.27'
$ws.Range('E20').Value = '  +0.59%  '
$ws.Range('E21').Value = '  +2.16%  '
$ws.Range('E22').Value = '  +2.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.52'
$ws.Range('E23').Value = '  +3.45%  '
$ws.Range('B24').Value = 'RenderToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.49'
$ws.Range('E24').Value = '  +10.72%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.22'
$ws.Range('E25').Value = '  +1.49%  '
$ws.Range('E26').Value = '  +7.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.39'
$ws.Range('E27').Value = '  +4.18%  '
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.57'
$ws.Range('E29').Value = '  +6.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.28'
$ws.Range('E30').Value = '  +11.83%  '
$ws.Range('E31').Value = '  +1.41%  '
$ws.Range('E32').Value = '  +0.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.57'
$ws.Range('E33').Value = '  +1.20%  '
$ws.Range('E34').Value = '  +3.36%  '
$ws.Range('E35').Value = '  +7.00%  '
$ws.Range('E36').Value = '  +3.03%  '
$ws.Range('E37').Value = '  +2.42%  '
$ws.Range('E38').Value = '  +11.12%  '
$ws.Range('E39').Value = '  +10.14%  '
$ws.Range('E40').Value = '  +3.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '50.48'
$ws.Range('E41').Value = '  +1.02%  '
$ws.Range('E42').Value = '  +1.59%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.309'
$ws.Range('E43').Value = '  +15.74%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '44.47'
$ws.Range('E44').Value = '  +15.63%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '391.52'
$ws.Range('E45').Value = '  +1.19%  '
$ws.Range('E46').Value = '  +3.77%  '
$ws.Range('D47').Value = '2.721.06'
$ws.Range('E47').Value = '  +1.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '133.02'
$ws.Range('E48').Value = '  +2.65%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.78'
$ws.Range('E49').Value = '  +11.51%  '
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('E51').Value = '  +7.69%  '
